# Dados apresentação 22-08.xlsx — refactor: Update MIGO test case and improve Excel data handling
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- A1: fix typo "Pedido Origen" -> "Pedido Origem" ---
$ws.Range("A1").Value = "Pedido Origem"

# --- Column B: was "Novo Pedido" (empty data) -> becomes "Preço", taking over
#     the values that used to live in column Q ("Preço": 6 / 7). Re-style the
#     header to match the normal blue headers (copy A1's format) instead of
#     the old green fill. ---
$ws.Range("B2").Value = $ws.Range("Q2").Value2
$ws.Range("B3").Value = $ws.Range("Q3").Value2
$ws.Range("A1").Copy()
$ws.Range("B1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B1").Value = "Preço"

# --- Column M: was "NOVO CONTRATO" (empty data) -> becomes "TP CONTRATO",
#     taking over the values that used to live in column R
#     ("TIP DE CONTRATO": ZDDR / ZDDR). Re-style the header like the other
#     blue headers instead of the old distinct fill. ---
$ws.Range("M2").Value = $ws.Range("R2").Value2
$ws.Range("M3").Value = $ws.Range("R3").Value2
$ws.Range("A1").Copy()
$ws.Range("M1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("M1").Value = "TP CONTRATO"

# --- Columns P,Q,R: old "Status"/"Preço"/"TIP DE CONTRATO" data (OK / 6-7 /
#     ZDDR) is cleared out (those values were relocated to B and M above),
#     and the headers are replaced with new blank test-case columns, all
#     sharing the same green header style as P1 already had. ---
$ws.Range("P2:R3").ClearContents()
$ws.Range("P1").Copy()
$ws.Range("Q1:R1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("P1").Value = "NV CONTRATO"
$ws.Range("Q1").Value = "NV PEDIDO"
$ws.Range("R1").Value = "MIGO"

# --- New column S: "MIRO" header, same green style as the other new headers ---
$ws.Range("P1").Copy()
$ws.Range("S1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("S1").Value = "MIRO"

# --- New placeholder cell L8: empty, underlined font, no fill ---
$ws.Range("L8").Font.Underline = $true

$excel.CutCopyMode = $false
